$wb = $excel.ActiveWorkbook

# Rename the existing sheet (was "AddMultipleCustomer")
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "addMultipleCustomerTest"

# Update header row text on sheet 1
$ws1.Range("A1").Value = "First Name"
$ws1.Range("B1").Value = "Last Name"
$ws1.Range("C1").Value = "Post Code"

# Add a new data row to sheet 1
$ws1.Range("A4").Value = "xyz"
$ws1.Range("B4").Value = "sample"
$ws1.Range("C4").NumberFormat = "@"
$ws1.Range("C4").Value = "00001"

# Update selection on sheet1 to match target (D9)
$ws1.Range("D9").Select()

# Add a second worksheet right after sheet 1
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "openAccountTest"
$ws2.Range("A1").Value = "Customer"
$ws2.Range("B1").Value = "Currency"
$ws2.Range("A1:B1").Style = $ws1.Range("A1").Style
$ws2.Range("A2").Select()
